{"js": "// Replace the \"150 words\" placeholder paragraph (under the \"Results\" heading)\n// with the full Results paragraph text, inserted as a sequence of runs that\n// mirror the author's original edit history.\nconst runs = [\"Our project aims to develop an android workout application with a focus on privacy and\", \" p\", \"ersonalized workout generation. It includes an Anonymous mode for privacy-conscious users and a normal mode with additional ChatGPT-powered personalized workout generation\", \" functionality\", \". Utilizing Android Studio and Java\", \" \u2013 the official IDE and the most used language for android application development\", \", we designed a user-friendly home page featuring a\", \"n \", \"Anonymous Browsing button to emphasize incognito usage. \", \"When one \", \"enters\", \" anonymously, he gains access to a prepopulated SQLite Database with various workouts to suit anyone\u2019s fitness needs\", \". \", \"On the other hand, when one logs in normally\", \", we implemented TiDB's online database\", \" f\", \"or secure user authentication\", \". \", \"Finally, u\", \"pon logging in, users can utilize various functionalities, including an AI workout generator powered by ChatGPT's API\", \". All they must do to receive a\", \" \", \"custom workout\", \" based on \", \"their\", \" preferences and goals\", \" is to answer a brief questionary.\"];\n\n// Locate the placeholder text left under the \"Results\" heading.\nconst search = context.document.body.search(\"150 words\", { matchCase: true });\nsearch.load(\"items\");\nawait context.sync();\n\nif (search.items.length === 0) {\n  throw new Error('Could not find placeholder text \"150 words\"');\n}\n\n// The first run replaces the placeholder text entirely; `insertText` hands\n// back a Range over exactly the text just inserted, so every later run can\n// be chained directly onto it (inserted immediately after the previous one)\n// without having to re-search the document.\nlet cursor = search.items[0].insertText(runs[0], Word.InsertLocation.replace);\nawait context.sync();\n\nfor (let i = 1; i < runs.length; i++) {\n  cursor = cursor.insertText(runs[i], Word.InsertLocation.after);\n  await context.sync();\n}\n", "ps1": "# Replace the \"150 words\" placeholder paragraph (under the \"Results\" heading)\n# with the full Results paragraph text, inserted as a sequence of runs that\n# mirror the author's original edit history.\n$runs = @(\n    \"Our project aims to develop an android workout application with a focus on privacy and\",\n    \" p\",\n    \"ersonalized workout generation. It includes an Anonymous mode for privacy-conscious users and a normal mode with additional ChatGPT-powered personalized workout generation\",\n    \" functionality\",\n    \". Utilizing Android Studio and Java\",\n    \" \u2013 the official IDE and the most used language for android application development\",\n    \", we designed a user-friendly home page featuring a\",\n    \"n \",\n    \"Anonymous Browsing button to emphasize incognito usage. \",\n    \"When one \",\n    \"enters\",\n    \" anonymously, he gains access to a prepopulated SQLite Database with various workouts to suit anyone\u2019s fitness needs\",\n    \". \",\n    \"On the other hand, when one logs in normally\",\n    \", we implemented TiDB's online database\",\n    \" f\",\n    \"or secure user authentication\",\n    \". \",\n    \"Finally, u\",\n    \"pon logging in, users can utilize various functionalities, including an AI workout generator powered by ChatGPT's API\",\n    \". All they must do to receive a\",\n    \" \",\n    \"custom workout\",\n    \" based on \",\n    \"their\",\n    \" preferences and goals\",\n    \" is to answer a brief questionary.\"\n)\n\n$d = $word.ActiveDocument\n$range = $d.Content\n$found = $range.Find.Execute(\"150 words\")\nif (-not $found) {\n    throw 'Could not find placeholder text \"150 words\"'\n}\n\n$range.Text = $runs[0]\nfor ($i = 1; $i -lt $runs.Length; $i++) {\n    $range.Collapse(0)  # wdCollapseEnd\n    $range.InsertAfter($runs[$i])\n}\n"}
